# ----------------------------------------------------------------------------
# Adds the 2022-Q4 quarterly snapshot:
#   1. Insert a new "2022-Q4" worksheet (cloned from "2022-Q3" so it keeps the
#      exact same formatting/styles), populated with the Q4-2022 fund holdings.
#   2. Prepend a new "2022-Q4" row to the "总计" (total) summary sheet, and
#      shift every existing row down by one.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# =============================================================================
# PART 1 -- "总计" (summary) sheet: add the 2022-Q4 row, push the rest down
# =============================================================================

$summary = $wb.Worksheets.Item("总计")

# Copy the index-column style (bold / bordered / centered) from the last
# existing data row down onto the new row 10 that we are about to populate.
$summary.Cells.Item(9, 1).Copy()
$summary.Cells.Item(10, 1).PasteSpecial(-4122)

# New summary table, top to bottom (row 2 = most recent quarter):
#   row, index, label,    count, value
$summaryRows = @(
    @(2, 0, "2022-Q4", 11, 0.65),
    @(3, 1, "2022-Q3", 10, 0.61),
    @(4, 2, "2022-Q2", 13, 1.25),
    @(5, 3, "2022-Q1", 14, 1.39),
    @(6, 4, "2021-Q4", 15, 1.26),
    @(7, 5, "2021-Q3", 9, 0.96),
    @(8, 6, "2021-Q2", 6, 0.58),
    @(9, 7, "2021-Q1", 18, 1.34),
    @(10, 8, "2020-Q4", 8, 0.31)
)

foreach ($r in $summaryRows) {
    $rowNum = $r[0]
    $summary.Cells.Item($rowNum, 1).Value = $r[1]
    $summary.Cells.Item($rowNum, 2).Value = $r[2]
    $summary.Cells.Item($rowNum, 3).Value = $r[3]
    $summary.Cells.Item($rowNum, 4).Value = $r[4]
}

# =============================================================================
# PART 2 -- new "2022-Q4" fund-holdings worksheet
# =============================================================================

$srcSheet = $wb.Worksheets.Item("2022-Q3")
$srcSheet.Copy($srcSheet)

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# The source sheet (old 2022-Q3) only had 10 data rows (rows 2-11); the new
# data has 11 data rows, so extend the index column's style onto row 12.
$q4.Cells.Item(11, 1).Copy()
$q4.Cells.Item(12, 1).PasteSpecial(-4122)

# Header row (row 1) is already correct on the cloned sheet (基金代码, 基金名称,
# 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名) -- no change needed.

# Data rows: index, code, name, size, position, pct, value, rank
$q4Rows = @(
    @(2, 0, "005775", "中加转型动力灵活配置混合A", "6.53", "53.27", "2.75", "0.1796", 8),
    @(3, 1, "630010", "华商价值精选混合", "4.37", "87.25", "3.33", "0.1455", 6),
    @(4, 2, "009242", "中加核心智造混合A", "1.95", "61.58", "3.78", "0.0737", 3),
    @(5, 3, "011815", "恒越优势精选混合", "2.64", "92.01", "2.30", "0.0607", 6),
    @(6, 4, "005776", "中加转型动力灵活配置混合C", "1.91", "53.27", "2.75", "0.0525", 8),
    @(7, 5, "012072", "中加喜利回报一年持有期混合C", "2.02", "36.36", "2.50", "0.0505", 6),
    @(8, 6, "012071", "中加喜利回报一年持有期混合A", "1.82", "36.36", "2.50", "0.0455", 6),
    @(9, 7, "630006", "华商产业升级混合", "0.86", "88.65", "3.40", "0.0292", 6),
    @(10, 8, "010130", "海富通惠增多策略一年定期开放灵活配置混合A", "0.51", "85.70", "1.99", "0.0101", 10),
    @(11, 9, "009243", "中加核心智造混合C", "0.11", "61.58", "3.78", "0.0042", 3),
    @(12, 10, "010131", "海富通惠增多策略一年定期开放灵活配置混合C", "0.00", "85.70", "1.99", "0.00", 10)
)

foreach ($r in $q4Rows) {
    $rowNum = $r[0]
    $q4.Cells.Item($rowNum, 1).Value = $r[1]
    $q4.Cells.Item($rowNum, 2).Value = "'" + $r[2]
    $q4.Cells.Item($rowNum, 3).Value = $r[3]
    $q4.Cells.Item($rowNum, 4).Value = "'" + $r[4]
    $q4.Cells.Item($rowNum, 5).Value = "'" + $r[5]
    $q4.Cells.Item($rowNum, 6).Value = "'" + $r[6]
    $q4.Cells.Item($rowNum, 7).Value = "'" + $r[7]
    $q4.Cells.Item($rowNum, 8).Value = $r[8]
}

# Last row's 持有市值(亿元) column (G12) is a genuine numeric 0, unlike every
# other row where that column is stored as text -- match the source data.
$q4.Cells.Item(12, 7).Value = 0
